$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "In Progress" formatting from F8 onto F7, and set F7's text accordingly.
$ws.Range("F8").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F7").Value = "In Progress"

$ws.Range("H7").Select() | Out-Null
